$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must remain stored as text,
# matching the original inline-string cells, not auto-converted numbers/percentages)
$updates = @{
    'D2' = '282.53'
    'E2' = '1.62%'
    'E3' = '4.21%'
    'D4' = '5.015'
    'E4' = '2.99%'
    'D5' = '0.06503'
    'E5' = '1.15%'
    'D6' = '7.215'
    'E6' = '3.15%'
    'D7' = '1.379'
    'E7' = '14.61%'
    'D8' = '0.9179'
    'E8' = '3.59%'
    'D9' = '0.1538'
    'E9' = '-0.73%'
    'D10' = '0.06288'
    'E10' = '22.51%'
    'D11' = '0.07564'
    'E11' = '0.72%'
    'D12' = '0.02826'
    'E12' = '-2.16%'
    'D13' = '0.08959'
    'E13' = '-0.16%'
    'D14' = '0.001590'
    'E14' = '0.56%'
    'D15' = '0.0006348'
    'E15' = '-0.25%'
    'D16' = '0.006101'
    'E16' = '-0.44%'
    'D17' = '3.445'
    'E17' = '-1.00%'
    'D18' = '3.357'
    'E18' = '1.59%'
    'D19' = '2.241'
    'E19' = '1.04%'
    'D21' = '0.1320'
    'E21' = '-2.24%'
    'D22' = '3.983'
    'E22' = '1.82%'
    'E23' = '2.90%'
    'D24' = '0.04428'
    'E24' = '0.21%'
    'D25' = '0.001182'
    'E25' = '0.51%'
    'D26' = '0.004456'
    'E26' = '14.91%'
    'E28' = '1.67%'
    'D29' = '0.0001617'
    'E29' = '-1.62%'
    'D40' = '0.04106'
    'E40' = '-0.43%'
    'D41' = '0.006658'
    'E41' = '-1.66%'
    'D42' = '0.1234'
    'E42' = '5.08%'
    'D43' = '0.002140'
    'E43' = '12.60%'
    'D44' = '0.01148'
    'E44' = '-3.06%'
    'D45' = '0.00005626'
    'E45' = '5.55%'
    'D46' = '1.950'
    'E46' = '16.08%'
    'E47' = '-0.09%'
}

foreach ($cellAddr in $updates.Keys) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cellAddr]
    $rng.Style = $origStyle
}
